$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, copying formatting (bold font, border, centered
# alignment) from the existing header cell G1 so the new column matches
# the style of the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add 0 values in H2:H8 for the new "Save" column (unstyled, like the
# other numeric data cells).
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
